# Atualização de bases das ligas, do dia: 24-02-2024 às 21:58
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Slovenia Prva Liga")

# --- Rows 80/81 swap their match data (id/date/div stay put) ---
$ws.Range("B80").Value = 5498503
$ws.Range("F80").Value = "FC Koper"
$ws.Range("G80").Value = "NS Mura"
$ws.Range("H80").Value = 1
$ws.Range("I80").Value = 2
$ws.Range("J80").Value = "A"
$ws.Range("K80").Value = 2.05
$ws.Range("L80").Value = 3.3
$ws.Range("M80").Value = 3.25
$ws.Range("N80").Value = 2
$ws.Range("O80").Value = 3.4
$ws.Range("P80").Value = 3.25
$ws.Range("Q80").Value = -0.5
$ws.Range("R80").Value = 2
$ws.Range("S80").Value = 1.8
$ws.Range("T80").Value = 2.5
$ws.Range("U80").Value = 1.825
$ws.Range("V80").Value = 1.975
$ws.Range("W80").Value = -1
$ws.Range("X80").Value = -1
$ws.Range("Y80").Value = 2.25
$ws.Range("Z80").Value = -1
$ws.Range("AA80").Value = 0.8
$ws.Range("AB80").Value = 0.825
$ws.Range("AC80").Value = -1

$ws.Range("B81").Value = 5495053
$ws.Range("F81").Value = "NK Radomlje"
$ws.Range("G81").Value = "NK Domzale"
$ws.Range("H81").Value = 1
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = "H"
$ws.Range("K81").Value = 2.55
$ws.Range("L81").Value = 3.1
$ws.Range("M81").Value = 2.55
$ws.Range("N81").Value = 3.75
$ws.Range("O81").Value = 3.4
$ws.Range("P81").Value = 1.833
$ws.Range("Q81").Value = 0.5
$ws.Range("R81").Value = 1.925
$ws.Range("S81").Value = 1.875
$ws.Range("T81").Value = 2.5
$ws.Range("U81").Value = 1.975
$ws.Range("V81").Value = 1.825
$ws.Range("W81").Value = 2.75
$ws.Range("X81").Value = -1
$ws.Range("Y81").Value = -1
$ws.Range("Z81").Value = 0.925
$ws.Range("AA81").Value = -1
$ws.Range("AB81").Value = -1
$ws.Range("AC81").Value = 0.825

# --- Rows 82/83 swap their match data (id/date/div stay put) ---
$ws.Range("B82").Value = 6816473
$ws.Range("F82").Value = "NK Bravo"
$ws.Range("G82").Value = "NK Rogaska"
$ws.Range("H82").Value = 2
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = "H"
$ws.Range("K82").Value = 1.8
$ws.Range("L82").Value = 3.5
$ws.Range("M82").Value = 4
$ws.Range("N82").Value = 2.05
$ws.Range("O82").Value = 3
$ws.Range("P82").Value = 3.75
$ws.Range("Q82").Value = -0.25
$ws.Range("R82").Value = 1.75
$ws.Range("S82").Value = 2.05
$ws.Range("T82").Value = 2.25
$ws.Range("U82").Value = 1.95
$ws.Range("V82").Value = 1.85
$ws.Range("W82").Value = 1.05
$ws.Range("X82").Value = -1
$ws.Range("Y82").Value = -1
$ws.Range("Z82").Value = 0.75
$ws.Range("AA82").Value = -1
$ws.Range("AB82").Value = -0.5
$ws.Range("AC82").Value = 0.425

$ws.Range("B83").Value = 6814327
$ws.Range("F83").Value = "NS Mura"
$ws.Range("G83").Value = "NK Domzale"
$ws.Range("H83").Value = 2
$ws.Range("I83").Value = 3
$ws.Range("J83").Value = "A"
$ws.Range("K83").Value = 2
$ws.Range("L83").Value = 3.3
$ws.Range("M83").Value = 3.4
$ws.Range("N83").Value = 1.909
$ws.Range("O83").Value = 3.4
$ws.Range("P83").Value = 3.75
$ws.Range("Q83").Value = -0.5
$ws.Range("R83").Value = 1.95
$ws.Range("S83").Value = 1.85
$ws.Range("T83").Value = 2.5
$ws.Range("U83").Value = 1.9
$ws.Range("V83").Value = 1.9
$ws.Range("W83").Value = -1
$ws.Range("X83").Value = -1
$ws.Range("Y83").Value = 2.75
$ws.Range("Z83").Value = -1
$ws.Range("AA83").Value = 0.8500000000000001
$ws.Range("AB83").Value = 0.8999999999999999
$ws.Range("AC83").Value = -1

# --- Row 188: match result now known (was unplayed) ---
$ws.Range("H188").Value = 0
$ws.Range("I188").Value = 4
$ws.Range("J188").Value = "A"
$ws.Range("N188").Value = 7
$ws.Range("O188").Value = 5
$ws.Range("P188").Value = 1.363
$ws.Range("R188").Value = 2
$ws.Range("S188").Value = 1.8
$ws.Range("T188").Value = 2.75
$ws.Range("U188").Value = 1.775
$ws.Range("V188").Value = 2.025
$ws.Range("W188").Value = -1
$ws.Range("X188").Value = -1
$ws.Range("Y188").Value = 0.363
$ws.Range("Z188").Value = -1
$ws.Range("AA188").Value = 0.8
$ws.Range("AB188").Value = 0.7749999999999999
$ws.Range("AC188").Value = -1

# --- Row 190: odds refresh ---
$ws.Range("N190").Value = 5.75
$ws.Range("P190").Value = 1.4

# --- Row 192: odds refresh ---
$ws.Range("N192").Value = 2.375
$ws.Range("P192").Value = 2.7
$ws.Range("R192").Value = 1.775
$ws.Range("S192").Value = 2.025
